$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vertical Sprite Mgmt")

# --- Update the core inputs/outputs for the vertical-sprite model ---
# cursl
$ws.Range("B1").Value = 15
# spypos
$ws.Range("B2").Value = 0
# mirrored (now TRUE)
$ws.Range("B4").Value = $true

# spyoff formula now subtracts 1 when mirrored
$ws.Range("B5").Formula = "=IF(B4,B3-(B1-B2)-1,B1-B2)"

# Update the pseudocode label for spyoff to match the new formula
$ws.Range("E4").Value = "spyoff = mir ? spysz - (cursl - spypos) - 1 : cursl - spypos"

# --- Add the new H/I/K/L helper table (rows 3-18) ---
# H column: 0..15 ; I column: FLOOR.MATH(H/8)
# K column: 15..0 ; L column: FLOOR.MATH(K/8)
for ($i = 0; $i -le 15; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 8).Value = $i
    $ws.Cells.Item($row, 9).Formula = "=_xlfn.FLOOR.MATH(H$row/8)"

    $kval = 15 - $i
    $ws.Cells.Item($row, 11).Value = $kval
    $ws.Cells.Item($row, 12).Formula = "=_xlfn.FLOOR.MATH(K$row/8)"
}

# --- Cosmetic: selection + column width, matching the saved state ---
$ws.Range("E2").Select()
$ws.Columns("E").AutoFit()
